$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 0.02201793066408441
$ws.Range("C3").Value = 0.02179891650764746
$ws.Range("D3").Value = 0.02179671840349664

# Row 4 - GradientBoostingRegressor
$ws.Range("C4").Value = 0.01257974491403423
$ws.Range("D4").Value = 0.01257974491403422

# Row 5 - AdaBoostRegressor
$ws.Range("B5").Value = 0.04654036700884014
$ws.Range("C5").Value = 0.04555487471020344
$ws.Range("D5").Value = 0.0408996534809167
